$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend / rewrite the "Kosten" data table (A2:C30) -----------------
# Column A: time in seconds, 5..145 step 5 (rows 2..30)
# Column B: "Beste Werte"    -> constant 35850 for every row
# Column C: "Optimale Werte" -> constant 15100 for every row
$row = 2
for ($zeit = 5; $zeit -le 145; $zeit += 5) {
    $ws.Cells.Item($row, 1).Value = $zeit
    $ws.Cells.Item($row, 2).Value = 35850
    $ws.Cells.Item($row, 3).Value = 15100
    $row++
}

# --- Update the chart series so they reference the new, larger range ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection()
for ($i = 1; $i -le $series.Count; $i++) {
    $ser = $series.Item($i)
    $f = $ser.Formula
    $f = $f.Replace('$A$2:$A$7', '$A$2:$A$30')
    $f = $f.Replace('$B$2:$B$7', '$B$2:$B$30')
    $f = $f.Replace('$C$2:$C$7', '$C$2:$C$30')
    $ser.Formula = $f
}
